$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.889.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.83%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.861.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.42%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'304.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.88%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.05%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5051"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.14%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3621"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.56%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07170"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.30%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8960"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.96%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.29%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07471"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.843.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.28%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'92.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.90%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.228"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.73%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.02%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008472"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.02%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'14.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.36%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.05%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'26.927.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.85%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.027"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.93%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'2.064.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.46%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -1.96%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.417"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.90%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'147.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.14%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.796"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.45%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'17.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.74%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -2.32%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'113.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.28%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.57%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.677"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.31%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.09260"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.47%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05086"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.80%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.13%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.7433"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.42%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.148"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.93%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.279"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.79%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'VeChain"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.01999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.69%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'TheSandbox"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.5606"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.47%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.11%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.63%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'118.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.46%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'6.476"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.65%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.508"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.95%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.1467"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.31%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4720"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.83%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.9998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.08%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'10.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.69%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.564"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.14%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'36.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.18%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -2.35%  "
$ws.Range("E51").Style = "Normal"
